$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename Sheet1 to "misc."
$ws.Name = "misc."

# Add two new header cells (K4, L4) for the new "other_indexes" / "commodity"
# columns, matching the style used by the rest of the header row.
$ws.Cells.Item(4, 11).Value = "other_indexes"
$ws.Cells.Item(4, 11).Style = "Heading 3"
$ws.Cells.Item(4, 12).Value = "commodity"
$ws.Cells.Item(4, 12).Style = "Heading 3"

# Add two new rows of data (row 11 and row 12).
# Write order matters for shared-string table ordering, so new unique
# strings are introduced in this sequence: flo_emis, gas, co2captured,
# co2, *ccs,*ccs-rf, coal,oil
$ws.Cells.Item(11, 2).Value = "flo_emis"
$ws.Cells.Item(11, 4).Value = "gas"
$ws.Cells.Item(11, 12).Value = "co2captured"
$ws.Cells.Item(11, 11).Value = "co2"
$ws.Cells.Item(11, 5).Value = "*ccs,*ccs-rf"
$ws.Cells.Item(12, 4).Value = "coal,oil"
$ws.Cells.Item(11, 8).Value = 0.95

$ws.Cells.Item(12, 2).Value = "flo_emis"
$ws.Cells.Item(12, 5).Value = "*ccs,*ccs-rf"
$ws.Cells.Item(12, 11).Value = "co2"
$ws.Cells.Item(12, 12).Value = "co2captured"
$ws.Cells.Item(12, 8).Value = 0.85

# Resize columns E and K to fit their new, wider content (matches Excel's
# own bestFit recalculation for the longest entries now present: "*ccs,*ccs-rf"
# in column E and "other_indexes" in column K).
$ws.Columns.Item(5).ColumnWidth = 9.498697916666666
$ws.Columns.Item(11).ColumnWidth = 11.166666666666666

# Update selection to D13
$ws.Range("D13").Select()
